# Update column C (Fitness) values in Sheet1 to reflect re-run results.
# Mapping is "row:newValue" pairs for all affected rows (2-178).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = "2:10470;3:10470;4:10385;5:10385;6:10385;7:10385;8:10385;9:9103;10:9103;11:9103;12:8943;13:8943;14:8532;15:8532;16:8532;17:8532;18:8532;19:8236;20:8236;21:8236;22:7741;23:7741;24:7741;25:7741;26:7741;27:7741;28:7741;29:7741;30:7741;31:7741;32:7741;33:7672;34:7672;35:7672;36:7672;37:7672;38:7672;39:7672;40:7672;41:7672;42:7672;43:7672;44:7672;45:7672;46:7672;47:7672;48:7672;49:7672;50:7672;51:7672;52:7672;53:7672;54:7672;55:7672;56:7672;57:7672;58:7672;59:7672;60:7672;61:7672;62:7672;63:7672;64:7672;65:7672;66:7672;67:7672;68:7672;69:7672;70:7672;71:7672;72:7672;73:7672;74:7639;75:7639;76:7639;77:7639;78:7639;79:7639;80:7639;81:7639;82:7639;83:7639;84:7639;85:7639;86:7639;87:7639;88:7639;89:7639;90:7639;91:7639;92:7639;93:7639;94:7639;95:7639;96:7639;97:7639;98:7639;99:7639;100:7639;101:7639;102:7312;103:7312;104:7312;105:7312;106:7312;107:7312;108:7312;109:7312;121:7295;122:7295;123:7295;124:7295;125:7295;126:7295;127:7295;128:7295;129:7295;130:7295;131:7295;132:7295;133:7295;134:7295;135:7295;136:7295;137:7295;138:7295;139:7295;140:7295;141:7295;142:7295;143:7295;144:7295;145:7295;146:7295;147:7295;148:7295;149:7295;150:7295;151:7295;152:7295;153:7295;154:7295;155:7295;156:7295;157:7295;158:7295;159:7295;160:7295;161:7295;162:7295;163:7295;164:7295;165:7295;166:7295;167:7295;168:7295;169:7295;170:7295;171:7295;172:7295;173:7295;174:7295;175:7295;176:7295;177:7295;178:7295"

foreach ($pair in $changes.Split(";")) {
    $parts = $pair.Split(":")
    $row = [int]$parts[0]
    $newValue = [int]$parts[1]
    $ws.Cells.Item($row, 3).Value = $newValue
}
